$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired content for rows 2-6 (row 1 header is unchanged)
$data = @(
    @("Georgia - Spain ✓: 0:4", "Spain", 74, "", 100, 1.25, "✓"),
    @("Cyprus - Austria ✓: 0:2", "Austria", 71, "", 93, 1.32, "✓"),
    @("CA Boca Juniors  - Club Atlético Tigre: 22:00", "CA Boca Juniors", 65, "", 94, 1.75, ""),
    @("Kazakhstan - Belgium : 1:1", "Belgium", 61, 75, 88, 1.18, ""),
    @("Racing Santander  - Granada CF: 2:2", "Racing Santander", 53, "", 81, 1.88, "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Rows 7-11 from the old table are no longer part of the data; remove them
$ws.Range("A7:G11").Delete()
